# Add "Inverter Alternatives" column (H) with a short list of alternative
# parts for the signal inverter, plus a couple of blank-but-formatted cells
# further down the sheet (H10:I10) picked up by a stray selection/format.
#
# The shared-string table order in the saved file mirrors the order in
# which distinct strings are first written, so we deliberately write H3
# before H1/H2/H4 to reproduce that order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "74AHC1G04SE-7 "
$ws.Range("H1").Value = "Inverter Alternatives:"
$ws.Range("H2").Value = "74LVC1G06SE-7"
$ws.Range("H4").Value = "TC7S04FU,LF "

# Match the "Text" number format (style index 1) already used by the other
# descriptive columns (A, D, E) on this sheet.
$ws.Range("H1:H4").NumberFormat = "@"

# A couple of stray formatted-but-empty cells a bit further down, same as
# in the authored workbook.
$ws.Range("H10:I10").NumberFormat = "@"

# Size the new column to fit its contents, like the other bestFit columns.
$ws.Columns("H:H").AutoFit() | Out-Null

# Leave the selection on the last-edited cell.
$ws.Range("H4").Select() | Out-Null
